# Regenerate the "K" (strikeouts) column (column G) in the save_data sheet.
# Commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals" -- the upstream data-generation script
# recomputed column G (header "K") for each outing row; this reproduces
# the resulting cell values in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 2
    13 = 0
    14 = 2
    15 = 2
    16 = 2
    17 = 1
    18 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
